# Update the "want to go" (想去人数) counts (column F) across the four
# worksheets to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 229
$ws.Range("F3").Value = 560
$ws.Range("F4").Value = 164
$ws.Range("F7").Value = 3127
$ws.Range("F8").Value = 2719
$ws.Range("F12").Value = 335
$ws.Range("F13").Value = 277
$ws.Range("F15").Value = 5583
$ws.Range("F16").Value = 608
$ws.Range("F17").Value = 1014
$ws.Range("F18").Value = 51
$ws.Range("F19").Value = 76
$ws.Range("F20").Value = 438
$ws.Range("F21").Value = 1199
$ws.Range("F23").Value = 104
$ws.Range("F24").Value = 323

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 12
$ws.Range("F9").Value = 48
$ws.Range("F13").Value = 627
$ws.Range("F25").Value = 4025
$ws.Range("F27").Value = 8
$ws.Range("F30").Value = 58

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2546
$ws.Range("F7").Value = 3
$ws.Range("F9").Value = 1420
$ws.Range("F10").Value = 397
$ws.Range("F11").Value = 110

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 2546
$ws.Range("F7").Value = 1420
$ws.Range("F8").Value = 397
$ws.Range("F9").Value = 110
$ws.Range("F11").Value = 12
$ws.Range("F12").Value = 229
$ws.Range("F13").Value = 560
$ws.Range("F14").Value = 164
$ws.Range("F16").Value = 3127
$ws.Range("F17").Value = 2719
$ws.Range("F20").Value = 335
$ws.Range("F22").Value = 48
$ws.Range("F23").Value = 277
$ws.Range("F25").Value = 5583
$ws.Range("F27").Value = 608
$ws.Range("F28").Value = 1014
$ws.Range("F29").Value = 627
$ws.Range("F30").Value = 51
$ws.Range("F31").Value = 76
$ws.Range("F32").Value = 438
$ws.Range("F40").Value = 1199
$ws.Range("F41").Value = 8
$ws.Range("F45").Value = 58
$ws.Range("F48").Value = 104
$ws.Range("F49").Value = 323
